$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$qSheet = $wb.Worksheets.Item("2022-Q2")

# The existing "2022-Q2" sheet becomes the new "2022-Q4" sheet; duplicate it first so an
# untouched "2022-Q2" sheet (keeping the original Q2 data/formatting) remains right after it.
$qSheet.Copy($null, $qSheet)
$newQ2Sheet = $wb.Worksheets.Item($qSheet.Index + 1)

# Rename: original sheet becomes "2022-Q4", the duplicate becomes the new "2022-Q2".
$qSheet.Name = "2022-Q4"
$newQ2Sheet.Name = "2022-Q2"

# The new quarter sheet's header row + A2 pick up the "总计" sheet's header styling
# (not the old Q2 sheet's styling that came along with the duplicate), so copy that
# formatting across instead of leaving the inherited one.
$totalSheet.Range("B1").Copy()
$qSheet.Range("B1:H1").PasteSpecial(-4122) # xlPasteFormats
$totalSheet.Range("A2").Copy()
$qSheet.Range("A2").PasteSpecial(-4122) # xlPasteFormats

# Overwrite the (now "2022-Q4") sheet's detail row with the new quarter's fund data.
# Columns B, D, E, F and G hold numeric-looking values that must stay text (e.g. the
# leading zero in the fund code), so round-trip them through a text formula + paste-as
# -values instead of a plain .Value assignment (which would auto-convert to a number and
# silently drop the leading zero / change the stored type).
$qSheet.Range("B2").Formula = '="015245"'
$qSheet.Range("D2").Formula = '="0.11"'
$qSheet.Range("E2").Formula = '="84.24"'
$qSheet.Range("F2").Formula = '="1.07"'
$qSheet.Range("G2").Formula = '="0.0012"'
$textRng = $qSheet.Range("B2:G2")
$textRng.Copy()
$textRng.PasteSpecial(-4163) # xlPasteValues

$qSheet.Range("C2").Value = "南华丰汇混合"
$qSheet.Range("H2").Value = 7

# Update the "总计" summary sheet: row 2 now reflects the newest quarter (2022-Q4), and a
# new row 3 is appended for the quarter that used to sit in row 2 (2022-Q2), matching the
# existing row's formatting.
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122) # xlPasteFormats
